$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Insert a new item row before current row 7 (CORASORE) for the
#    new "BRUFEN" item, copying the row-7 (post header) formatting.
# ------------------------------------------------------------------
$ws.Rows.Item(7).Insert()
$ws.Range("A8:Q8").Copy()
$ws.Range("A7:Q7").PasteSpecial(-4122)

$ws.Range("A7").Value = 1
$ws.Range("C7").Value = "BRUFEN 100MG/5ML SYRUP 150ML"
$ws.Range("H7").Value = "1:0"
$ws.Range("L7").Value = "1"
$ws.Range("N7").Value = "44.00"
$ws.Range("P7").Value = "44.0000"
$ws.Range("Q7").Value = "1:0"

# Renumber the rows that followed (CORASORE .. NANAZOXID), now rows 8-11
$ws.Range("A8").Value = 2
$ws.Range("A9").Value = 3
$ws.Range("A10").Value = 4
$ws.Range("A11").Value = 5

# ------------------------------------------------------------------
# 2) Insert a new item row before the "سرنجات 3 سم" row (now row 12)
#    for the new "NEOZOLID" item, copying that row's formatting.
# ------------------------------------------------------------------
$ws.Rows.Item(12).Insert()
$ws.Range("A13:Q13").Copy()
$ws.Range("A12:Q12").PasteSpecial(-4122)

$ws.Range("A12").Value = 6
$ws.Range("C12").Value = "NEOZOLID 250 MG 6 TABS."
$ws.Range("H12").Value = "0:0"
$ws.Range("L12").Value = "1"
$ws.Range("N12").Value = "36.50"
$ws.Range("P12").Value = "36.5000"
$ws.Range("Q12").Value = "1:0"

# "سرنجات 3 سم" row shifted from 12 to 13; renumber it
$ws.Range("A13").Value = 7

# ------------------------------------------------------------------
# 3) Update the totals cell (now row 14) and the generated-at
#    timestamp footer (now row 15).
# ------------------------------------------------------------------
$ws.Range("P14").Value = 171

$ws.Range("A15").Value = "Friday, 11 July, 2025 5:02 PM"

Write-Output "edit complete"
